$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# Header updates: "Ultima actualizacion" timestamp and "Total filas" count
$ws1.Cells.Item(2,1).Value = "Última actualización: 09:38:09"
$ws1.Cells.Item(3,1).Value = "Total filas: 150"

$ws2.Cells.Item(2,1).Value = "Última actualización: 09:38:09"
$ws2.Cells.Item(3,1).Value = "Total filas: 20"

$ws3.Cells.Item(2,1).Value = "Última actualización: 09:38:09"
$ws3.Cells.Item(3,1).Value = "Total filas: 29"

# Sheet "LP1912": rows 42-43 swapped; rows 117-155 rescheduled/added (38 rows touched)
$ws1.Cells.Item(42,1).Value = "06:57:30"
$ws1.Cells.Item(42,2).Value = "06:58"
$ws1.Cells.Item(42,3).Value = "14_ABASTO"
$ws1.Cells.Item(42,4).Value = 1
$ws1.Cells.Item(42,5).Value = "LP1912"
$ws1.Cells.Item(43,1).Value = "05:18:56"
$ws1.Cells.Item(43,2).Value = "06:58"
$ws1.Cells.Item(43,3).Value = "10_OLMOS"
$ws1.Cells.Item(43,4).Value = 100
$ws1.Cells.Item(43,5).Value = "LP1912"
$ws1.Cells.Item(117,1).Value = "08:39:44"
$ws1.Cells.Item(117,2).Value = "09:34"
$ws1.Cells.Item(117,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(117,4).Value = 55
$ws1.Cells.Item(117,5).Value = "LP1912"
$ws1.Cells.Item(118,1).Value = "08:39:44"
$ws1.Cells.Item(118,2).Value = "09:34"
$ws1.Cells.Item(118,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(118,4).Value = 55
$ws1.Cells.Item(118,5).Value = "LP1912"
$ws1.Cells.Item(122,1).Value = "09:38:09"
$ws1.Cells.Item(122,2).Value = "09:41"
$ws1.Cells.Item(122,3).Value = "14_ABASTO"
$ws1.Cells.Item(122,4).Value = 3
$ws1.Cells.Item(122,5).Value = "LP1912"
$ws1.Cells.Item(123,1).Value = "09:38:09"
$ws1.Cells.Item(123,2).Value = "09:41"
$ws1.Cells.Item(123,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(123,4).Value = 3
$ws1.Cells.Item(123,5).Value = "LP1912"
$ws1.Cells.Item(124,1).Value = "07:47:32"
$ws1.Cells.Item(124,2).Value = "09:42"
$ws1.Cells.Item(124,3).Value = "215C_EL PATO"
$ws1.Cells.Item(124,4).Value = 115
$ws1.Cells.Item(124,5).Value = "LP1912"
$ws1.Cells.Item(125,1).Value = "07:47:32"
$ws1.Cells.Item(125,2).Value = "09:43"
$ws1.Cells.Item(125,3).Value = "14_ABASTO"
$ws1.Cells.Item(125,4).Value = 116
$ws1.Cells.Item(125,5).Value = "LP1912"
$ws1.Cells.Item(126,1).Value = "08:57:13"
$ws1.Cells.Item(126,2).Value = "09:44"
$ws1.Cells.Item(126,3).Value = "14_ABASTO"
$ws1.Cells.Item(126,4).Value = 47
$ws1.Cells.Item(126,5).Value = "LP1912"
$ws1.Cells.Item(127,1).Value = "09:38:09"
$ws1.Cells.Item(127,2).Value = "09:47"
$ws1.Cells.Item(127,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(127,4).Value = 9
$ws1.Cells.Item(127,5).Value = "LP1912"
$ws1.Cells.Item(128,1).Value = "08:49:51"
$ws1.Cells.Item(128,2).Value = "09:52"
$ws1.Cells.Item(128,3).Value = "15_ABASTO"
$ws1.Cells.Item(128,4).Value = 63
$ws1.Cells.Item(128,5).Value = "LP1912"
$ws1.Cells.Item(129,1).Value = "08:49:51"
$ws1.Cells.Item(129,2).Value = "09:53"
$ws1.Cells.Item(129,3).Value = "10_OLMOS"
$ws1.Cells.Item(129,4).Value = 64
$ws1.Cells.Item(129,5).Value = "LP1912"
$ws1.Cells.Item(130,1).Value = "09:38:09"
$ws1.Cells.Item(130,2).Value = "09:59"
$ws1.Cells.Item(130,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(130,4).Value = 21
$ws1.Cells.Item(130,5).Value = "LP1912"
$ws1.Cells.Item(131,1).Value = "09:38:09"
$ws1.Cells.Item(131,2).Value = "10:04"
$ws1.Cells.Item(131,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(131,4).Value = 26
$ws1.Cells.Item(131,5).Value = "LP1912"
$ws1.Cells.Item(132,1).Value = "09:38:09"
$ws1.Cells.Item(132,2).Value = "10:05"
$ws1.Cells.Item(132,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(132,4).Value = 27
$ws1.Cells.Item(132,5).Value = "LP1912"
$ws1.Cells.Item(133,1).Value = "08:39:44"
$ws1.Cells.Item(133,2).Value = "10:06"
$ws1.Cells.Item(133,3).Value = "10_OLMOS"
$ws1.Cells.Item(133,4).Value = 87
$ws1.Cells.Item(133,5).Value = "LP1912"
$ws1.Cells.Item(134,1).Value = "08:21:50"
$ws1.Cells.Item(134,2).Value = "10:10"
$ws1.Cells.Item(134,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(134,4).Value = 109
$ws1.Cells.Item(134,5).Value = "LP1912"
$ws1.Cells.Item(135,1).Value = "08:57:13"
$ws1.Cells.Item(135,2).Value = "10:11"
$ws1.Cells.Item(135,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(135,4).Value = 74
$ws1.Cells.Item(135,5).Value = "LP1912"
$ws1.Cells.Item(136,1).Value = "08:21:50"
$ws1.Cells.Item(136,2).Value = "10:12"
$ws1.Cells.Item(136,3).Value = "15_ABASTO"
$ws1.Cells.Item(136,4).Value = 111
$ws1.Cells.Item(136,5).Value = "LP1912"
$ws1.Cells.Item(137,1).Value = "09:38:09"
$ws1.Cells.Item(137,2).Value = "10:13"
$ws1.Cells.Item(137,3).Value = "10_OLMOS"
$ws1.Cells.Item(137,4).Value = 35
$ws1.Cells.Item(137,5).Value = "LP1912"
$ws1.Cells.Item(138,1).Value = "08:49:51"
$ws1.Cells.Item(138,2).Value = "10:20"
$ws1.Cells.Item(138,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(138,4).Value = 91
$ws1.Cells.Item(138,5).Value = "LP1912"
$ws1.Cells.Item(139,1).Value = "08:39:44"
$ws1.Cells.Item(139,2).Value = "10:21"
$ws1.Cells.Item(139,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(139,4).Value = 102
$ws1.Cells.Item(139,5).Value = "LP1912"
$ws1.Cells.Item(140,1).Value = "08:39:44"
$ws1.Cells.Item(140,2).Value = "10:22"
$ws1.Cells.Item(140,3).Value = "17_ROMERO"
$ws1.Cells.Item(140,4).Value = 103
$ws1.Cells.Item(140,5).Value = "LP1912"
$ws1.Cells.Item(141,1).Value = "09:38:09"
$ws1.Cells.Item(141,2).Value = "10:24"
$ws1.Cells.Item(141,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(141,4).Value = 46
$ws1.Cells.Item(141,5).Value = "LP1912"
$ws1.Cells.Item(142,1).Value = "08:39:44"
$ws1.Cells.Item(142,2).Value = "10:26"
$ws1.Cells.Item(142,3).Value = "215A_EL PATO"
$ws1.Cells.Item(142,4).Value = 107
$ws1.Cells.Item(142,5).Value = "LP1912"
$ws1.Cells.Item(143,1).Value = "08:57:13"
$ws1.Cells.Item(143,2).Value = "10:27"
$ws1.Cells.Item(143,3).Value = "215A_EL PATO"
$ws1.Cells.Item(143,4).Value = 90
$ws1.Cells.Item(143,5).Value = "LP1912"
$ws1.Cells.Item(144,1).Value = "08:49:51"
$ws1.Cells.Item(144,2).Value = "10:41"
$ws1.Cells.Item(144,3).Value = "17_ROMERO"
$ws1.Cells.Item(144,4).Value = 112
$ws1.Cells.Item(144,5).Value = "LP1912"
$ws1.Cells.Item(145,1).Value = "08:57:13"
$ws1.Cells.Item(145,2).Value = "10:42"
$ws1.Cells.Item(145,3).Value = "17_ROMERO"
$ws1.Cells.Item(145,4).Value = 105
$ws1.Cells.Item(145,5).Value = "LP1912"
$ws1.Cells.Item(146,1).Value = "08:49:51"
$ws1.Cells.Item(146,2).Value = "10:43"
$ws1.Cells.Item(146,3).Value = "14_ABASTO"
$ws1.Cells.Item(146,4).Value = 114
$ws1.Cells.Item(146,5).Value = "LP1912"
$ws1.Cells.Item(147,1).Value = "08:57:13"
$ws1.Cells.Item(147,2).Value = "10:44"
$ws1.Cells.Item(147,3).Value = "14_ABASTO"
$ws1.Cells.Item(147,4).Value = 107
$ws1.Cells.Item(147,5).Value = "LP1912"
$ws1.Cells.Item(148,1).Value = "09:38:09"
$ws1.Cells.Item(148,2).Value = "10:58"
$ws1.Cells.Item(148,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(148,4).Value = 80
$ws1.Cells.Item(148,5).Value = "LP1912"
$ws1.Cells.Item(149,1).Value = "09:38:09"
$ws1.Cells.Item(149,2).Value = "11:02"
$ws1.Cells.Item(149,3).Value = "215C_EL PATO"
$ws1.Cells.Item(149,4).Value = 84
$ws1.Cells.Item(149,5).Value = "LP1912"
$ws1.Cells.Item(150,1).Value = "09:38:09"
$ws1.Cells.Item(150,2).Value = "11:07"
$ws1.Cells.Item(150,3).Value = "16_P MOR-167 Y 521"
$ws1.Cells.Item(150,4).Value = 89
$ws1.Cells.Item(150,5).Value = "LP1912"
$ws1.Cells.Item(151,1).Value = "09:38:09"
$ws1.Cells.Item(151,2).Value = "11:20"
$ws1.Cells.Item(151,3).Value = "86_EST CHICA-ESC AGRARIA"
$ws1.Cells.Item(151,4).Value = 102
$ws1.Cells.Item(151,5).Value = "LP1912"
$ws1.Cells.Item(152,1).Value = "09:38:09"
$ws1.Cells.Item(152,2).Value = "11:21"
$ws1.Cells.Item(152,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(152,4).Value = 103
$ws1.Cells.Item(152,5).Value = "LP1912"
$ws1.Cells.Item(153,1).Value = "09:38:09"
$ws1.Cells.Item(153,2).Value = "11:27"
$ws1.Cells.Item(153,3).Value = "225_C ROCA-H SUR"
$ws1.Cells.Item(153,4).Value = 109
$ws1.Cells.Item(153,5).Value = "LP1912"
$ws1.Cells.Item(154,1).Value = "09:38:09"
$ws1.Cells.Item(154,2).Value = "11:32"
$ws1.Cells.Item(154,3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(154,4).Value = 114
$ws1.Cells.Item(154,5).Value = "LP1912"
$ws1.Cells.Item(155,1).Value = "09:38:09"
$ws1.Cells.Item(155,2).Value = "11:36"
$ws1.Cells.Item(155,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(155,4).Value = 118
$ws1.Cells.Item(155,5).Value = "LP1912"

# Sheet "LP1912-215": new row appended (1 row(s) touched)
$ws2.Cells.Item(25,1).Value = "09:38:09"
$ws2.Cells.Item(25,2).Value = "11:02"
$ws2.Cells.Item(25,3).Value = "215C_EL PATO"
$ws2.Cells.Item(25,4).Value = 84
$ws2.Cells.Item(25,5).Value = "LP1912"

# Sheet "6203-6173": new row appended (1 row(s) touched)
$ws3.Cells.Item(34,1).Value = "09:38:09"
$ws3.Cells.Item(34,2).Value = "11:14"
$ws3.Cells.Item(34,3).Value = "215C_LA PLATA"
$ws3.Cells.Item(34,4).Value = 96
$ws3.Cells.Item(34,5).Value = "L6203"

